$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("229").Delete()
